$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-24 Wednesday" "2024-07-25 Thursday"

Replace-Text "140×6=" "676×9="
Replace-Text "345×6=" "423×3="
Replace-Text "253×2=" "816×3="
Replace-Text "394×2=" "483×7="
Replace-Text "393×2=" "133×8="
Replace-Text "895×7=" "720×2="
Replace-Text "522×4=" "521×6="
Replace-Text "690×3=" "400×2="
Replace-Text "353×2=" "115×3="
Replace-Text "160×5=" "571×9="
Replace-Text "645×4=" "640×4="
Replace-Text "906×7=" "812×5="
Replace-Text "531×2=" "221×9="
Replace-Text "487×4=" "586×6="
Replace-Text "726×5=" "797×2="
Replace-Text "318×5=" "316×3="
Replace-Text "869×4=" "806×4="
Replace-Text "216×9=" "857×2="
Replace-Text "891×8=" "598×5="
Replace-Text "376×3=" "944×2="
Replace-Text "517×2=" "116×2="
Replace-Text "358×2=" "218×8="
Replace-Text "537×3=" "752×9="
Replace-Text "263×2=" "939×7="
Replace-Text "690×7=" "666×3="
